$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.616372666666667
$ws.Range("H2").Value = 13.849118
$ws.Range("I2").Value = 0.0183283362562958
$ws.Range("J2").Value = 0.01832833625629581
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 15.95540881354755
$ws.Range("R2").Value = 143.598679321928
$ws.Range("S2").Value = 0.000180378977561489
$ws.Range("T2").Value = 0.000180378977561489

$ws.Range("G3").Value = 4.616372666666667
$ws.Range("H3").Value = 13.849118
$ws.Range("I3").Value = 0.0183283362562958
$ws.Range("J3").Value = 0.01832833625629581
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 1392.298682567403
$ws.Range("R3").Value = 12530.68814310663
$ws.Range("S3").Value = 0.01574020545361865
$ws.Range("T3").Value = 0.01574020545361865

$ws.Range("G4").Value = 4.616372666666667
$ws.Range("H4").Value = 13.849118
$ws.Range("I4").Value = 0.0183283362562958
$ws.Range("J4").Value = 0.01832833625629581
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 212.9775055310425
$ws.Range("R4").Value = 1916.797549779382
$ws.Range("S4").Value = 0.002407751825115672
$ws.Range("T4").Value = 0.002407751825115672

$ws.Range("I5").Value = 0.943783113604627
$ws.Range("J5").Value = 0.9437831136046271
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 821.5936895915485
$ws.Range("R5").Value = 7394.343206323936
$ws.Range("S5").Value = 0.009288275307221299
$ws.Range("T5").Value = 0.009288275307221299

$ws.Range("I6").Value = 0.943783113604627
$ws.Range("J6").Value = 0.9437831136046271
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.8105121983829771
$ws.Range("T6").Value = 0.8105121983829769

$ws.Range("I7").Value = 0.943783113604627
$ws.Range("J7").Value = 0.9437831136046271
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.1239826399144288
$ws.Range("T7").Value = 0.1239826399144288

$ws.Range("I8").Value = 0.03788855013907712
$ws.Range("J8").Value = 0.03788855013907712
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 32.98320689712978
$ws.Range("R8").Value = 296.848862074168
$ws.Range("S8").Value = 0.0003728815228947118
$ws.Range("T8").Value = 0.0003728815228947118

$ws.Range("I9").Value = 0.03788855013907712
$ws.Range("J9").Value = 0.03788855013907712
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.03253833600548169
$ws.Range("T9").Value = 0.03253833600548169

$ws.Range("I10").Value = 0.03788855013907712
$ws.Range("J10").Value = 0.03788855013907712
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.004977332610700727
$ws.Range("T10").Value = 0.004977332610700726
